$wb = $excel.ActiveWorkbook
$props = @("Width","Height","UsableWidth","UsableHeight","Left","Top","WindowState")
foreach ($p in $props) {
   try {
      Write-Output "$p => $($excel.$p)"
   } catch { Write-Output "$p ERR $_" }
}
